$d = $word.ActiveDocument

# --- Add the new "Tomorrow's weather" paragraph after the last paragraph ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last

# Type the paragraph text, plus a temporary trailing sentinel character. The
# sentinel keeps the true end-of-paragraph position away from the (buggy)
# "one before the paragraph mark" offset while we create the bookmark, so the
# bookmark lands exactly where we want instead of snapping to the start of
# the document.
$newPara.Range.Text = "Tomorrow" + [char]0x2019 + "s weather: Probably even colder.~"

$withSentinel = $d.Paragraphs.Last
$bookmarkPos = $withSentinel.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)

# "_GoBack" is the singleton "last edit" bookmark; adding it here both places
# it on the new paragraph and implicitly removes it from wherever it used to
# be (the "Yesterday's weather" paragraph).
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the temporary sentinel character now that the bookmark is anchored.
$withSentinel2 = $d.Paragraphs.Last
$sentinelRange = $d.Range($withSentinel2.Range.End - 2, $withSentinel2.Range.End - 1)
$sentinelRange.Delete()
